$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "aydan"
$ws.Range("B5").Value = "jk"

$ws.Range("B5").Select()
